$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Iterator Traceability" data table (student 2 / iteration trace), appended
# below the existing Bridge/Factory/Iterator sections (previous last row was 24).
# Header row first...
$ws.Cells.Item(25, 1).Value = "Java Code"
$ws.Cells.Item(25, 2).Value = "C# Code"
$ws.Cells.Item(25, 3).Value = "Test Case(s)"

# ...then the first data row's Java/C# cells...
$ws.Cells.Item(26, 1).Value = "J17"
$ws.Cells.Item(26, 2).Value = "C17"

# ...then the rest of the Java Code column...
$javaVals = @("J18", "J19", "J20", "J21", "J22", "J23", "J24", "J25", "J26", "J27", "J28")
$r = 27
foreach ($v in $javaVals) {
    $ws.Cells.Item($r, 1).Value = $v
    $r++
}

# ...then the rest of the C# Code column...
$csVals = @("C18", "C19", "C20", "C21", "C22", "C23", "C24", "C25", "C26", "C27", "C28")
$r = 27
foreach ($v in $csVals) {
    $ws.Cells.Item($r, 2).Value = $v
    $r++
}

# ...then the Test Case(s) column for every data row.
$testVals = @(
    "T14, T15, T16, T17, T18, T19",
    "T14, T15, T16, T17, T18, T19",
    "T14, T15, T16, T17",
    "T14, T15, T16, T17",
    "T14, T15, T16, T17",
    "T16, T17",
    "T14, T15, T16, T17",
    "T16, T17",
    "T14, T15, T16, T17, T18, T19",
    "T19",
    "T18",
    "T14, T15, T16, T17"
)
$r = 26
foreach ($v in $testVals) {
    $ws.Cells.Item($r, 3).Value = $v
    $r++
}

$startRow = 25
$endRow = 37

# Match the bordered "data row" formatting used throughout the sheet (e.g. row 2)
# by copying its format onto the freshly written rows.
$ws.Range("A2:C2").Copy()
$ws.Range("A" + $startRow + ":C" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the viewport/selection where the author left it after adding the new rows.
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("B44").Select()
